# "Add key and certificates"
#
# Data refresh of the Expenses sheet:
#   - Row 3: new id/user_id guids, amount 10 -> 11, description cleared, timestamps bumped
#   - Row 4: new id guid, user_id guid now matches row 3's, amount 10 -> 12,
#            description cleared, timestamps bumped
#   - Rows 5-7: three new expense rows appended, same user_id/category pattern,
#            each with its own amount + timestamp pair and a blank description
#
# Existing rows 3 & 4 already carry the date-number style (from F3/G3) and the
# default/general style (from A1), so most of this script only needs to set
# .Value - styles are picked up implicitly. New rows 5-7 don't exist yet, so
# their date cells need the date style copied over explicitly (Range.Style
# assignment isn't reliable here, so we drive it via Copy/PasteSpecial like a
# real user would with "Paste Special > Formats").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blankCellStyleSource = $ws.Range("A1")   # default/general style, no number format
$dateCellStyleSource  = $ws.Range("F3")   # existing date-formatted style (s="1")

function Set-BlankDescription($cell) {
    # Typing a bare apostrophe is how Excel records an explicit empty text
    # value (quote-prefixed blank) instead of clearing the cell outright.
    $cell.Value = "'"
    $blankCellStyleSource.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats - strip the quote-prefix style back to normal
}

function Apply-DateStyle($range) {
    $dateCellStyleSource.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

function Set-ExpenseRow($r, $id, $userId, $amount, $category, $timestamp) {
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $userId
    $ws.Cells.Item($r, 3).Value = $amount
    Set-BlankDescription $ws.Cells.Item($r, 4)
    $ws.Cells.Item($r, 5).Value = $category
    $ws.Cells.Item($r, 6).Value = $timestamp
    $ws.Cells.Item($r, 7).Value = $timestamp
    Apply-DateStyle ($ws.Range($ws.Cells.Item($r, 6), $ws.Cells.Item($r, 7)))
}

Set-ExpenseRow 3 "71f07443-f6a6-4fef-a0a3-93b494e4e3af" "b68d09e7-df5d-4d13-bee3-ac455e4a33e9" 11 "Food" 45197.28138627315
Set-ExpenseRow 4 "6574d207-b6e8-4ad4-b1e8-f46ca7f4f0f3" "b68d09e7-df5d-4d13-bee3-ac455e4a33e9" 12 "Food" 45197.28142216435
Set-ExpenseRow 5 "20da1323-1d1c-4ec5-9394-2ca5f663702b" "b68d09e7-df5d-4d13-bee3-ac455e4a33e9" 45 "Food" 45197.28147271991
Set-ExpenseRow 6 "82e368a2-975e-4b5f-a7cd-3e5982d35bde" "b68d09e7-df5d-4d13-bee3-ac455e4a33e9" 34 "Food" 45197.28188016204
Set-ExpenseRow 7 "a013816c-0f7b-44e7-9c76-97cef7f30f1e" "b68d09e7-df5d-4d13-bee3-ac455e4a33e9" 56 "Food" 45197.28191664352
